$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 1).Value = "山子高科"
$ws.Cells.Item(2, 2).Value = "上海建工"
$ws.Cells.Item(2, 3).Value = "利欧股份"
$ws.Cells.Item(3, 1).Value = "上海建工"
$ws.Cells.Item(3, 2).Value = "山子高科"
$ws.Cells.Item(3, 3).Value = "岩山科技"
$ws.Cells.Item(4, 1).Value = "首开股份"
$ws.Cells.Item(4, 2).Value = "利欧股份"
$ws.Cells.Item(4, 3).Value = "卧龙电驱"
$ws.Cells.Item(5, 1).Value = "三花智控"
$ws.Cells.Item(5, 2).Value = "岩山科技"
$ws.Cells.Item(5, 3).Value = "山子高科"
$ws.Cells.Item(6, 1).Value = "岩山科技"
$ws.Cells.Item(6, 2).Value = "首开股份"
$ws.Cells.Item(6, 3).Value = "领益智造"
$ws.Cells.Item(7, 1).Value = "利欧股份"
$ws.Cells.Item(7, 2).Value = "供销大集"
$ws.Cells.Item(7, 3).Value = "青山纸业"
$ws.Cells.Item(8, 1).Value = "中科曙光"
$ws.Cells.Item(8, 2).Value = "华胜天成"
$ws.Cells.Item(8, 3).Value = "中国电影"
$ws.Cells.Item(9, 1).Value = "卧龙电驱"
$ws.Cells.Item(9, 2).Value = "中科曙光"
$ws.Cells.Item(9, 3).Value = "吉视传媒"
$ws.Cells.Item(10, 1).Value = "青山纸业"
$ws.Cells.Item(10, 2).Value = "三花智控"
$ws.Cells.Item(10, 3).Value = "先导智能"
$ws.Cells.Item(11, 1).Value = "胜宏科技"
$ws.Cells.Item(11, 2).Value = "卧龙电驱"
$ws.Cells.Item(11, 3).Value = "天际股份"
$ws.Cells.Item(12, 1).Value = "省广集团"
$ws.Cells.Item(12, 2).Value = "万向钱潮"
$ws.Cells.Item(12, 3).Value = "拓维信息"
$ws.Cells.Item(13, 1).Value = "汉威科技"
$ws.Cells.Item(13, 2).Value = "怡 亚 通"
$ws.Cells.Item(13, 3).Value = "胜宏科技"
$ws.Cells.Item(14, 1).Value = "华胜天成"
$ws.Cells.Item(14, 2).Value = "天赐材料"
$ws.Cells.Item(14, 3).Value = "东方财富"
$ws.Cells.Item(15, 1).Value = "供销大集"
$ws.Cells.Item(15, 2).Value = "大洋电机"
$ws.Cells.Item(15, 3).Value = "指南针"
$ws.Cells.Item(16, 1).Value = "中国电影"
$ws.Cells.Item(16, 2).Value = "均胜电子"
$ws.Cells.Item(16, 3).Value = "北方稀土"
$ws.Cells.Item(17, 1).Value = "领益智造"
$ws.Cells.Item(17, 2).Value = "青山纸业"
$ws.Cells.Item(17, 3).Value = "三花智控"
$ws.Cells.Item(18, 1).Value = "天赐材料"
$ws.Cells.Item(18, 2).Value = "省广集团"
$ws.Cells.Item(18, 3).Value = "中科曙光"
$ws.Cells.Item(19, 1).Value = "东方财富"
$ws.Cells.Item(19, 2).Value = "东方财富"
$ws.Cells.Item(19, 3).Value = "卧龙新能"
$ws.Cells.Item(20, 1).Value = "吉视传媒"
$ws.Cells.Item(20, 2).Value = "汉威科技"
$ws.Cells.Item(20, 3).Value = "华胜天成"
$ws.Cells.Item(21, 1).Value = "上海贝岭"
$ws.Cells.Item(21, 2).Value = "上海贝岭"
$ws.Cells.Item(21, 3).Value = "二六三"
